$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 64: date 2026-01-06 (serial 46028) and associated counts
$ws.Range("A64").Value = 46028
$ws.Range("B64").Value = 5599
$ws.Range("C64").Value = 4344
$ws.Range("D64").Value = 4028
$ws.Range("E64").Value = 236
$ws.Range("F64").Value = 45
$ws.Range("G64").Value = 30
$ws.Range("H64").Value = 4
$ws.Range("I64").Value = 1

# Row 65: date 2026-01-07 (serial 46029) and associated counts
$ws.Range("A65").Value = 46029
$ws.Range("B65").Value = 5590
$ws.Range("C65").Value = 4371
$ws.Range("D65").Value = 4065
$ws.Range("E65").Value = 227
$ws.Range("F65").Value = 40
$ws.Range("G65").Value = 35
$ws.Range("H65").Value = 3
$ws.Range("I65").Value = 1

# Update the active selection to match the newly entered last row
$ws.Range("A65:I65").Select()
